# Tab field filter repository - added basic insert.
# The "Field name" column (column C on Sheet1) is being removed entirely,
# which shifts the subsequent "Variant column name"/"Relation"/"Default value"/
# "Options" columns one position to the left (D->C, E->D, F->E) and drops
# the now-unused "Field name" / "Field 1xx" shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the whole "Field name" column (column C) - this shifts D:G left to C:F.
$ws.Columns.Item(3).Delete()

# Move the active selection to C2 (was F14 before the edit).
[void]$ws.Range("C2").Select()
